$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 417544.9
$ws.Range("J17").Value = 426407.56
$ws.Range("L17").Value = 1279222.68
$ws.Range("N17").Value = -1279558.68
$ws.Range("H38").Value = 1176.1904
$ws.Range("I38").Value = 407.2143
$ws.Range("K38").Value = 1221.6429
$ws.Range("M38").Value = -849.6428999999998
$ws.Range("H39").Value = 185.8
$ws.Range("I39").Value = 170
$ws.Range("J39").Value = 249
$ws.Range("K39").Value = 510
$ws.Range("L39").Value = 747
$ws.Range("M39").Value = -214
$ws.Range("N39").Value = -1339
$ws.Range("H53").Value = 413.53845
$ws.Range("I53").Value = 325.47058
$ws.Range("J53").Value = 579.8889
$ws.Range("K53").Value = 325.47058
$ws.Range("L53").Value = 579.8889
$ws.Range("M53").Value = 311.52942
$ws.Range("N53").Value = -1853.8889
$ws.Range("H80").Value = 2703.9375
$ws.Range("I80").Value = 527.5
$ws.Range("K80").Value = 1582.5
$ws.Range("M80").Value = -584.5
$ws.Range("H82").Value = 4599.8
$ws.Range("I82").Value = 5499.75
$ws.Range("J82").Value = 1000
$ws.Range("K82").Value = 16499.25
$ws.Range("L82").Value = 3000
$ws.Range("M82").Value = -16093.25
$ws.Range("N82").Value = -3812
$ws.Range("H83").Value = 2703.9375
$ws.Range("I83").Value = 527.5
$ws.Range("K83").Value = 4747.5
$ws.Range("M83").Value = 244.5
$ws.Range("H85").Value = 4599.8
$ws.Range("I85").Value = 5499.75
$ws.Range("J85").Value = 1000
$ws.Range("K85").Value = 16499.25
$ws.Range("L85").Value = 3000
$ws.Range("M85").Value = -15095.25
$ws.Range("N85").Value = -5808
$ws.Range("H86").Value = 1827.0435
$ws.Range("J86").Value = 1568.7778
$ws.Range("L86").Value = 1568.7778
$ws.Range("N86").Value = -3814.7778
$ws.Range("H88").Value = 1524.8334
$ws.Range("J88").Value = 1907
$ws.Range("L88").Value = 1907
$ws.Range("N88").Value = -2719
$ws.Range("H89").Value = 1827.0435
$ws.Range("J89").Value = 1568.7778
$ws.Range("L89").Value = 7843.889
$ws.Range("N89").Value = -19075.889
$ws.Range("H91").Value = 1524.8334
$ws.Range("J91").Value = 1907
$ws.Range("L91").Value = 1907
$ws.Range("N91").Value = -4715
$ws.Range("H100").Value = 4303.4707
$ws.Range("I100").Value = 4462.8335
$ws.Range("J100").Value = 3921
$ws.Range("K100").Value = 4462.8335
$ws.Range("L100").Value = 3921
$ws.Range("M100").Value = -3921.8335
$ws.Range("N100").Value = -5003
$ws.Range("H132").Value = 372929.8
$ws.Range("I132").Value = 587816.9399999999
$ws.Range("K132").Value = 1763450.82
$ws.Range("M132").Value = -1760920.82
$ws.Range("H135").Value = 7421.5654
$ws.Range("I135").Value = 3053.2
$ws.Range("K135").Value = 27478.8
$ws.Range("M135").Value = -24943.8
$ws.Range("H138").Value = 3604.7737
$ws.Range("J138").Value = 5267.514
$ws.Range("L138").Value = 15802.542
$ws.Range("N138").Value = -26082.542
$ws.Range("H141").Value = 1366.6364
$ws.Range("I141").Value = 1428
$ws.Range("J141").Value = 1203
$ws.Range("K141").Value = 4284
$ws.Range("L141").Value = 3609
$ws.Range("M141").Value = 896
$ws.Range("N141").Value = -13969

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 324.75
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = ""
$ws.Range("H88").Value = 14899.2
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 14899.2
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 14899.2
$ws.Range("N88").Value = -15711.2
$ws.Range("M88").Value = ""
$ws.Range("H91").Value = 14899.2
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 14899.2
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 14899.2
$ws.Range("N91").Value = -17707.2
$ws.Range("M91").Value = ""
$ws.Range("H132").Value = 1206805
$ws.Range("I132").Value = 1517414
$ws.Range("J132").Value = 171441.67
$ws.Range("K132").Value = 4552242
$ws.Range("L132").Value = 514325.01
$ws.Range("M132").Value = -4549712
$ws.Range("N132").Value = -519385.01

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6098.773
$ws.Range("J86").Value = 9316.666999999999
$ws.Range("L86").Value = 9316.666999999999
$ws.Range("N86").Value = -11562.667
$ws.Range("H89").Value = 6098.773
$ws.Range("J89").Value = 9316.666999999999
$ws.Range("L89").Value = 46583.335
$ws.Range("N89").Value = -57815.335
$ws.Range("H94").Value = 6664.3335
$ws.Range("I94").Value = 3393.25
$ws.Range("J94").Value = 7598.9287
$ws.Range("K94").Value = 3393.25
$ws.Range("L94").Value = 7598.9287
$ws.Range("M94").Value = -2942.25
$ws.Range("N94").Value = -8500.9287
$ws.Range("H134").Value = 1358726
$ws.Range("I134").Value = 1619561.8
$ws.Range("J134").Value = 11074.5
$ws.Range("K134").Value = 4858685.4
$ws.Range("L134").Value = 33223.5
$ws.Range("M134").Value = -4856150.4
$ws.Range("N134").Value = -38293.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6407.5835
$ws.Range("I80").Value = 4893.625
$ws.Range("J80").Value = 9435.5
$ws.Range("K80").Value = 4893.625
$ws.Range("L80").Value = 9435.5
$ws.Range("M80").Value = -3895.625
$ws.Range("N80").Value = -11431.5
$ws.Range("H83").Value = 6407.5835
$ws.Range("I83").Value = 4893.625
$ws.Range("J83").Value = 9435.5
$ws.Range("K83").Value = 24468.125
$ws.Range("L83").Value = 47177.5
$ws.Range("M83").Value = -19476.125
$ws.Range("N83").Value = -57161.5
$ws.Range("H102").Value = 5104.1836
$ws.Range("I102").Value = 4057.077
$ws.Range("K102").Value = 4057.077
$ws.Range("M102").Value = -2435.077
$ws.Range("H132").Value = 30307078
$ws.Range("I132").Value = 45458296
$ws.Range("K132").Value = 136374888
$ws.Range("M132").Value = -136372358

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5877
$ws.Range("I40").Value = 5314.8335
$ws.Range("K40").Value = 5314.8335
$ws.Range("M40").Value = -5178.8335
$ws.Range("H82").Value = 3574.5417
$ws.Range("I82").Value = 1171.3334
$ws.Range("J82").Value = 7579.8887
$ws.Range("K82").Value = 1171.3334
$ws.Range("L82").Value = 7579.8887
$ws.Range("M82").Value = -810.3334
$ws.Range("N82").Value = -8301.8887
$ws.Range("H85").Value = 3574.5417
$ws.Range("I85").Value = 1171.3334
$ws.Range("J85").Value = 7579.8887
$ws.Range("K85").Value = 1171.3334
$ws.Range("L85").Value = 7579.8887
$ws.Range("M85").Value = 76.66660000000002
$ws.Range("N85").Value = -10075.8887
$ws.Range("H132").Value = 4374.378
$ws.Range("I132").Value = 3432.8572
$ws.Range("J132").Value = 7669.7
$ws.Range("K132").Value = 10298.5716
$ws.Range("L132").Value = 23009.1
$ws.Range("M132").Value = -7768.571599999999
$ws.Range("N132").Value = -28069.1
$ws.Range("H136").Value = 45462172
$ws.Range("I136").Value = 62508436
$ws.Range("K136").Value = 187525308
$ws.Range("M136").Value = -187522758

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1654.1072
$ws.Range("I81").Value = 1571.4584
$ws.Range("J81").Value = 2150
$ws.Range("K81").Value = 3142.9168
$ws.Range("L81").Value = 4300
$ws.Range("M81").Value = -2081.9168
$ws.Range("N81").Value = -6422
$ws.Range("H84").Value = 1654.1072
$ws.Range("I84").Value = 1571.4584
$ws.Range("J84").Value = 2150
$ws.Range("K84").Value = 15714.584
$ws.Range("L84").Value = 21500
$ws.Range("M84").Value = -10410.584
$ws.Range("N84").Value = -32108
$ws.Range("H100").Value = 2306.4546
$ws.Range("I100").Value = 1232.8235
$ws.Range("K100").Value = 2465.647
$ws.Range("M100").Value = -1924.647
$ws.Range("H132").Value = 5102.4595
$ws.Range("I132").Value = 4283.154
$ws.Range("J132").Value = 7039
$ws.Range("K132").Value = 12849.462
$ws.Range("L132").Value = 21117
$ws.Range("M132").Value = -10319.462
$ws.Range("N132").Value = -26177
$ws.Range("H136").Value = 8202935
$ws.Range("I136").Value = 11905789
$ws.Range("K136").Value = 35717367
$ws.Range("M136").Value = -35714817

